# Track preorders: append newly pre-ordered Hot Wheels to the list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: S.No, Model Name, Series
$data = @(
    @(280, "'18 Honda Civic Type R", "Mainlines"),
    @(281, "Custom Datsun 240Z", "Mainlines"),
    @(282, "'71 Nissan Skyline H/T 2000GT-R", "Silver Series National Icons"),
    @(283, "Nissan Silvia (S15)", "Mainlines"),
    @(284, "Alpine A110", "Mainlines"),
    @(285, "'06 Honda Civic Si", "The Hot Ones"),
    @(286, "Toyota AE86 Sprinter Trueno", "Ultra Hots"),
    @(287, "Corvette Stingray ('76)", "Mainlines"),
    @(288, "'18 Ford Mustang RTR Spec 5", "Silver Series Mustang 60 Years")
)

$startRow = 281
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $item = $data[$i]

    $numCell = $ws.Cells.Item($row, 1)
    $numCell.Value = $item[0]

    $nameCell = $ws.Cells.Item($row, 2)
    $nameText = [string]$item[1]
    if ($nameText.StartsWith("'")) {
        # A leading apostrophe is normally consumed by Excel as a text-qualifier
        # prefix rather than stored as literal text. Doubling it up front makes
        # the first one act as the (discarded) prefix marker and keeps the
        # second as real text, then clearing Style removes the quotePrefix flag.
        $nameCell.Value = "'" + $nameText
    } else {
        $nameCell.Value = $nameText
    }
    $nameCell.Style = "Normal"

    $seriesCell = $ws.Cells.Item($row, 3)
    $seriesText = [string]$item[2]
    if ($seriesText.StartsWith("'")) {
        $seriesCell.Value = "'" + $seriesText
    } else {
        $seriesCell.Value = $seriesText
    }
    $seriesCell.Style = "Normal"
}
